$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 6.5
$ws.Range("C3").Value = 8.5

# Update selection to D7
$ws.Range("D7").Select()

$wb.Save()
